$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
  'Aafje Thuiszorg Huizen Zorghotels (Stichting)',
  'Accolade (Stichting)',
  'ActiVite (Stichting)',
  'Adelante Groep (Stichting)',
  'Alkcare (Stichting)',
  'Altrecht (Stichting)',
  'Amaliazorg (Stichting)',
  'Amarant Groep (Stichting)',
  'Amaris Zorggroep (Stichting)',
  'Amerpoort (Stichting)',
  'Archipel (Stichting)',
  'Argos Zorggroep (Stichting)',
  'ASVZ (Stichting)',
  'Atlant Zorggroep (Stichting)',
  'Attent',
  'Aveleijn (Stichting)',
  'Avoord Zorg en Wonen (Stichting)',
  'AxionContinu Groep (Stichting)',
  'Azora (Stichting)',
  'Baalderborg Groep (Stichting)',
  'Bartholomeus Gasthuis (Stichting)',
  'Basalt revalidatie',
  'Bethanië (Stichting)',
  'Beweging 3.0 (Stichting)',
  'BrabantZorg',
  'Cardia (Stichting)',
  'Careyn (Stichting)',
  'Carinova (Stichting)',
  'Carint-Reggeland Groep (Stichting)',
  'Charim (Zorggroep)',
  'Cicero Zorggroep (Stichting)',
  'Combinatie Jeugdzorg',
  'Cordaan Groep (Stichting)',
  'Cosis',
  'Curamare (Stichting)',
  'De Hoop ggz (Stichting)',
  'De Zijlen (Stichting) (Ilmarinen)',
  'DFZS De Forensische Zorgspecialisten',
  'Dichterbij (Stichting)',
  'Dienstencentrum Oud Burgeren Gasthuis (OBG) (Stichting)',
  'Dimence Groep (Stichting)',
  'Driegasthuizengroep',
  'Driestroom (Stichting)',
  'DSV (Stichting)',
  'DZN B.V. (directe zorg nijmegen)',
  'Elver (Stichting)',
  'Emergis (Stichting)',
  'Espria (stichting ...) (met onderdelen Trans, Meander, Evean, GGZ Drenthe, icare)',
  'Fier Fryslan (Stichting)',
  'Fundis (Stichting) (beheren meerdere zorginstellingen, zoals Welthuis)',
  'GGNet (Stichting)',
  'GGz Breburg Groep (Stichting)',
  'GGz Centraal (Stichting)',
  'GGZ Delfland - Geestelijke Gezondheidszorg Delfland (Stichting)',
  'GGZ Friesland',
  'GGZ inGeest (Stichting)',
  'GGZ Westelijk Noord-Brabant / GGZ-WNB (Stichting)',
  'GGzE (Stichting) GGz Eindhoven',
  'Hartekamp Groep (Stichting)',
  'Heliomare (Stichting)',
  'Het Gasthuis Millingen aan de Rijn (Sint Jan De Deo) (Stichting)',
  'Het Laar (Stichting)',
  'Het Parkhuis (Stichting)',
  'HilverZorg (Stichting)',
  'Hoogstraat Revalidatie (Stichting)',
  'Innoforte (Stichting)',
  'Interzorg Noord Nederland (Stichting)',
  'Ipse de Bruggen (Stichting)',
  'Kalorama (Stichting)',
  'Karakter (Stichting)',
  'Kennemerhart',
  'Klein Geluk',
  'Klimmendaal (Stichting)',
  'Koninklijke Visio, expertisecentrum voor slechtziende en blinde mensen (Stichting)',
  'Koperhorst (Stichting)',
  'KwadrantGroep (Stichting)',
  'Land van Horne (Stichting voor Verpleeg-, Verzorgings- en Woonfaciliteiten ...)',
  'Leger des Heils Welzijns- en Gezondheidszorg (Stichting)',
  'Lelie Zorggroep (Stichting)',
  'Lentis incl. Dignis',
  'Levvel (noord holland)',
  'Liante (Stichting)',
  'Libra Revalidatie & Audiologie (Stichting)',
  'Liemerije (Stichting)',
  'Lister (Stichting)',
  'Magenta',
  'Marente (Stichting)',
  'Mediant, Stichting voor Geestelijke Gezondheidszorg Oost- en Midden Twente',
  'MET-GGZ (Limburg)',
  'Middin (Stichting)',
  'Mijzo Schakelring/Eikendonk (Stichting) -> fusie Mijzo (met Riethorst + Volckaert)',
  'Mondriaan (Stichting)',
  'NNCZ (Noord Nederlandse Coöperatie van Zorgorganisaties)',
  'Noorderboog (Stichting)',
  'Noorderbreedte',
  'Novadic-Kentron (Stichting)',
  'Odion (Stichting)',
  'Omega (Groep, Zwolle!)',
  'Omring (Stichting)',
  'Opbouw (Stichting) incl. Prinsenstichting',
  'Ouderenzorg Oudewater, De Wulverhorst',
  'Pameijer (Stichting)',
  'Pantein (Stichting)',
  'Park Zuiderhout (Stichting)',
  'Parnassia Groep B.V. Incl. Parnassia haaglanden en noord holland, Antes, Brijder, Youz, etc.',
  'Pergamijn (Stichting)',
  'Philadelphia Zorg (Stichting)',
  'Pieter Raat Stichting',
  'Pieter van Foreest (Stichting)',
  'Pleyade (Stichting)',
  'Pluryn Hoenderloo Groep (Stichting)',
  'Prisma (Stichting)',
  'Pro Persona',
  'Pro Senectute (Stichting)',
  'Proteion Groep (Stichting)',
  'PSW (Stichting)',
  'QuaRijn (Stichting)',
  'R.K. Zorgcentrum Roomburgh (Stichting)',
  'Raffy-Leystroom',
  'Reade',
  'Residentie Molenwijck',
  'Respect Zorggroep (Stichting)',
  'Revalidatie Friesland',
  'Revant (Stichting)',
  'Rijnhoven (Stichting)',
  'Rivierduinen',
  'Robert Coppes Stichting',
  'S&L Zorg (Stichting)',
  'Saffier - De Residentie (Stichting Zorginstelling ...))',
  'Salem Verpleeghuis (Stichting)',
  'Savant, Organisatie voor Zorg (Stichting)',
  'Schärwachter B.V.',
  'SEIN Stichting Epilepsie Instellingen Nederland',
  'Severinusstichting',
  'SGL (Stichting)',
  'sHeerenLoo Zorggroep (Stichting)',
  'Sint Anna Boxmeer (Stichting)',
  'Sint Jacob (Stichting)',
  'Siza (Stichting)',
  'Solis (Stichting)',
  'SOVAK (Stichting)',
  'Sterk Huis (Stichting) (West Brabant was voorheen Juzt)',
  'Surplus (Stichting en surplus zorg)',
  'SVRZ (Stichting Voor Regionale Zorgverlening)',
  'Tactus Verslavingszorg (Stichting)',
  'tanteLouise (Stichting)',
  'Terwille verslavingszorg (Stichting)',
  'Thebe (Zorggroep west en midden Brabant, incl. Ruitersbos)',
  'Topaz (Stichting)',
  'Trajectum (Stichting)',
  'Treant Zorggroep (Stichting)',
  'Valkenhof (Stichting)',
  'Van Neynselstichting (Stichting)',
  'Vecht & Ijssel (Stichting)',
  'Verpleeghuis Bergweide (Stichting)',
  'Viersprong (Netherlands institute for personality disorders)',
  'Vilente (Stichting)',
  'Vincent van Gogh (Stichting)',
  'Vitalis',
  'ViVa! Zorggroep (Stichting)',
  'Vogellanden, Centrum voor Revalidatie (Stichting)',
  'Waardeburgh (Stichting)',
  'Warande (Stichting)',
  'WelThuis BV',
  'Werkt voor Ouderen (Stichting) (WVO Zorg)',
  'Wever (Stichting)',
  'WIJdezorg (Stichting)',
  'WilgaerdenLeekerweide Groep (Stichting)',
  'Wonen en Zorg Purmerend (Stichting) (SWZP)',
  'Woon en zorgcentrum de Merwelanden, stichting',
  'Woonzorg Samen (Stichting)',
  'Youke',
  'Yulius (Stichting)',
  'Zellingen (Stichting Zorgbeheer De ...)',
  'ZGR (Zorggroep Raalte (Stichting))',
  'Zonnehuisgroep Noord (Stichting)',
  'Zorgaccent (Stichting)',
  'Zorgbalans',
  'Zorgboog (Stichting)',
  'Zorgcentrum het Bildt (Beukelaar) (Stichting)',
  'Zorgfederatie Oldenzaal (Stichting)',
  'Zorggroep Amsterdam Oost (ZGAO) (Stichting)',
  'Zorggroep Apeldoorn en omstreken (Stichting)',
  'Zorggroep Elde Maasduinen (Maasduinen Zorg => gefuseerd uit GD HvB gestapt)',
  'Zorggroep Ena (Stichting)',
  'Zorggroep Groningen (Stichting)',
  'Zorggroep Sint Maarten (Stichting)',
  'Zorggroep Sirjon',
  'Zorggroep Tangenborgh (Stichting)',
  'Zorggroep Tellus (Stichting)',
  'Zorgpartners Midden-Holland (Stichting)',
  'ZorgSpectrum (Stichting)',
  'Zorgspectrum Het Zand',
  'Zozijn Beheer (Stichting)',
  'ZuidOostZorg (Stichting)',
  'ZZG Zorggroep (Stichting)'
)

$statuses = @(
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Vastgesteld',
  'Voorlopig',
  'Voorlopig'
)

for ($i = 0; $i -lt $names.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $names[$i]
  $ws.Cells.Item($row, 2).Value = $statuses[$i]
}

$ws.Range("A1").Select()